$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-09-04 Thursday" "2025-09-05 Friday"
Replace-Text "71×47=3337" "67×53=3551"
Replace-Text "95×22=2090" "40×72=2880"
Replace-Text "52×95=4940" "60×53=3180"
Replace-Text "74×22=1628" "91×53=4823"
Replace-Text "29×65=1885" "95×98=9310"
Replace-Text "47×12=564" "84×36=3024"
Replace-Text "83×96=7968" "65×85=5525"
Replace-Text "23×84=1932" "40×55=2200"
Replace-Text "28×83=2324" "75×43=3225"
Replace-Text "28×55=1540" "79×77=6083"
Replace-Text "78×94=7332" "25×85=2125"
Replace-Text "65×19=1235" "89×84=7476"
Replace-Text "96×82=7872" "12×85=1020"
Replace-Text "14×77=1078" "18×88=1584"
Replace-Text "44×88=3872" "17×91=1547"
Replace-Text "38×90=3420" "24×36=864"
Replace-Text "93×71=6603" "94×25=2350"
Replace-Text "71×48=3408" "62×71=4402"
Replace-Text "82×13=1066" "39×26=1014"
Replace-Text "64×60=3840" "17×41=697"
Replace-Text "45×17=765" "24×55=1320"
Replace-Text "83×28=2324" "61×11=671"
Replace-Text "17×22=374" "34×81=2754"
Replace-Text "63×45=2835" "64×69=4416"
Replace-Text "18×83=1494" "78×50=3900"
